$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "SubArraysWithEqual1sAnd0s"
$ws.Range("A20").Value = "Subarrays with equal 1s and 0s"

$ws.Range("A15").Select()
